$d = $word.ActiveDocument

# --- 1. Remove the first two empty paragraphs that follow the table ---
# (done first, while paragraph indices still reflect the untouched document;
#  inserting the new table row below shifts/otherwise disturbs that indexing)
$p = $d.Paragraphs.Item(22)
$p.Range.Delete()
$p = $d.Paragraphs.Item(22)
$p.Range.Delete()

$t = $d.Tables.Item(1)

# --- 2. Insert a new "NAL Unrated" row right after the header row ---
$newRow = $t.Rows.Add($t.Rows.Item(2))
$t.Cell(2, 1).Range.Text = "NAL Unrated"
$t.Cell(2, 2).Range.Text = "NAL Unrated games are either not yet rated or do not need to be rated. Please be aware that even though the game is not rated, it might still contain content found in higher rated games."

# --- 3. Update the "NAL 13" row description text ---
$range = $d.Content
$range.Find.Execute(" or alcohol can also be present.", $true, $false, $false, $false, $false, $true, 1, $false, ", alcohol or drugs can also be present.", 2) | Out-Null

$range = $d.Content
$range.Find.Execute(" or really mild drug use..", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null
